# Add the "Level 1" space-bar bug test case row to the "Level 1 Testing" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Level 1 Testing")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Level 1 > After pressing space button"
$ws.Range("C2").Value = "It should attack or move forward"
$ws.Range("D2").Value = "Charater stops at one point and doesn" + [char]8217 + "t move further once pressed spacebar key"

$ws.Range("D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 86.4

$ws.Activate() | Out-Null
$ws.Range("A2:XFD2").Select() | Out-Null
